$d = $word.ActiveDocument

function Find-ParaIndex($text) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs($i)
        if ($p.Range.Text -eq "$text`r") {
            return $i
        }
    }
    return -1
}

function Set-GreenShading($rng) {
    $rng.Shading.Texture = 0
    $rng.Shading.ForegroundPatternColor = -16777216
    $rng.Shading.BackgroundPatternColor = 5296274
}

# --- 1. Shade "...Buscar por ingredientes" and add two new shaded list items after it ---
$idx = Find-ParaIndex("...Buscar por ingredientes")
$pBuscar = $d.Paragraphs($idx)
Set-GreenShading($pBuscar.Range)

$pBuscar.Range.InsertParagraphAfter()
$pSubirNew = $d.Paragraphs($idx + 1)
$pSubirNew.Range.Text = "...Subir receta"
Set-GreenShading($pSubirNew.Range)

$pSubirNew.Range.InsertParagraphAfter()
$pValorarNew = $d.Paragraphs($idx + 2)
$pValorarNew.Range.Text = "...Valorar receta"
Set-GreenShading($pValorarNew.Range)

# --- 2. Remove the old (now duplicated) "...Subir receta" and "...Valorar receta" paragraphs ---
# (the old "Subir receta" paragraph also carries the stray _GoBack bookmark, which is removed
#  along with it)
$oldSubirIdx = Find-ParaIndex("...Subir receta")
while ($d.Paragraphs($oldSubirIdx).Range.Shading.BackgroundPatternColor -eq 5296274) {
    # skip the freshly-inserted shaded copy, look further down
    $oldSubirIdx = -1
    for ($i = $idx + 3; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -eq "...Subir receta`r") {
            $oldSubirIdx = $i
            break
        }
    }
    break
}
$d.Paragraphs($oldSubirIdx).Range.Delete()
# "...Valorar receta" is now shifted into the same slot
$d.Paragraphs($oldSubirIdx).Range.Delete()

# --- 3. Renumber the h.i3rtj9ajnsk5 bookmark's paragraph: add bookmarkEnd right after it, and
#         move the _GoBack bookmark onto the page-break paragraph that follows ---
$pageBreakIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $xml = $p.Range.WordOpenXML
    if ($xml -like "*w:br w:type=`"page`"*") {
        $pageBreakIdx = $i
        break
    }
}
$pPageBreak = $d.Paragraphs($pageBreakIdx)
$rBreak = $d.Range($pPageBreak.Range.Start, $pPageBreak.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $rBreak)

Write-Output "Done. Final paragraph count: $($d.Paragraphs.Count)"
